# Update "Generate Report for Handback" timestamps in the handback-status workbook.
$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Overview!G4 - "Latest HO Xliff Generate Date" for the 2dd397bb... file.
# This value was originally shared with de-de!H4 (same underlying text),
# so both cells must be updated together to the same new value.
$wsOverview.Range("G4").Value = "2016-08-21 06:51:59"
$wsDeDe.Range("H4").Value = "2016-08-21 06:51:59"

# zh-cn sheet, row for 2dd397bb... file
# H4 - Correspond Handoff Datetime
$wsZhCn.Range("H4").Value = "2016-08-21 06:51:54"
# K4 - Correspond Handback DateTime
$wsZhCn.Range("K4").Value = "2016-08-21 06:52:14"

# de-de sheet, row for 2dd397bb... file
# K4 - Correspond Handback DateTime
$wsDeDe.Range("K4").Value = "2016-08-21 06:52:21"
